# "As classes foram devidamente separadas"
# Row 3 ("2525.0"/"2412.4"/"14", stored as text) becomes proper numeric
# values, and three more rows of data are appended below it: two fully
# numeric rows (4 and 5) and a final row (6) that - like the original
# row 3 - is stored as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: convert from text to real numbers.
$ws.Cells.Item(3, 1).Value = 2525
$ws.Cells.Item(3, 2).Value = 2412.4
$ws.Cells.Item(3, 3).Value = 14

# Row 4: new numeric row.
$ws.Cells.Item(4, 1).Value = 5000
$ws.Cells.Item(4, 2).Value = 5000
$ws.Cells.Item(4, 3).Value = 20

# Row 5: new numeric row.
$ws.Cells.Item(5, 1).Value = 2000
$ws.Cells.Item(5, 2).Value = 2000
$ws.Cells.Item(5, 3).Value = 20

# Row 6: new row, stored as text (mirrors the original row-3 style).
# Force text format first so Excel doesn't auto-coerce these
# numeric-looking strings back into numbers.
$ws.Range("A6:C6").NumberFormat = "@"
$ws.Cells.Item(6, 1).Value = "8000.0"
$ws.Cells.Item(6, 2).Value = "8000.0"
$ws.Cells.Item(6, 3).Value = "50"
